# Qatar Stars League workbook update (17-05-2024 13:59)
# The underlying data source re-sorted several match rows; this manifests in the
# worksheet as groups of rows exchanging their full record (columns B:AB) while
# each row keeps its own sequential index in column A.
#
# Below, each "cycle" lists worksheet row numbers in rotation order: the NEW
# content of cycle[i] is the OLD content of cycle[i+1] (wrapping around).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cycles = @(
    ,@(24, 25)
    ,@(27, 28)
    ,@(37, 38)
    ,@(42, 43)
    ,@(45, 46)
    ,@(68, 69)
    ,@(94, 95)
    ,@(96, 97)
    ,@(110, 111)
    ,@(122, 127, 123, 126, 124)
    ,@(128, 129)
)

foreach ($cycle in $cycles) {
    # Capture the original (pre-edit) values of column B:AB for every row in the cycle.
    $originals = @{}
    foreach ($r in $cycle) {
        $rng = $ws.Range("B$r`:AB$r")
        $originals[$r] = $rng.Value()
    }

    # Write back so that each row receives the values that belonged to the
    # next row in the cycle (rotation), leaving column A untouched.
    $count = $cycle.Length
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $cycle[$i]
        $srcRow = $cycle[($i + 1) % $count]
        $destRng = $ws.Range("B$destRow`:AB$destRow")
        $destRng.Value = $originals[$srcRow]
    }
}
